$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the duplicate (unused) chart-helper defined names left over from the
# chart's data range being re-picked.
$null = $wb.Names.Item("_xlchart.v1.2").Delete()
$null = $wb.Names.Item("_xlchart.v1.3").Delete()

# New header labels for the "increase" summary block (D18:F18), bold style
# like the other section headers (D15:E15 "Std" / "Relative std").
$ws.Range("D15:E15").Copy()
$ws.Range("D18:F18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D18").Value = "Mean increase"
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = "Median increase"

# Formulas for the increase values (D19 / F19)
$ws.Range("D19").Formula = "=((E3 / 114.202998) * 100) - 100"
$ws.Range("D19").ClearFormats()
$ws.Range("F19").Formula = "=((E10 / 113.658804) * 100) - 100"
$ws.Range("F19").ClearFormats()

$ws.Calculate()

# Final selection left on the sheet after the edit
$null = $ws.Range("G23").Select()
